# Auto-generated Excel COM-interop edit script
# Applies per-cell numeric updates (market price refresh) across 8 Leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 206.57143
$ws.Range("I4").Value = 232.66667
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 232.66667
$ws.Range("L4").Value = 50
$ws.Range("M4").Value = -118.66667
$ws.Range("N4").Value = -278
$ws.Range("H11").Value = 538.75
$ws.Range("I11").Value = 538.75
$ws.Range("K11").Value = 538.75
$ws.Range("M11").Value = -398.75
$ws.Range("H20").Value = 1925.8572
$ws.Range("H21").Value = 17
$ws.Range("I21").Value = 17
$ws.Range("K21").Value = 17
$ws.Range("M21").Value = 451
$ws.Range("H23").Value = 17
$ws.Range("I23").Value = 17
$ws.Range("K23").Value = 17
$ws.Range("M23").Value = 217
$ws.Range("H26").Value = 12000
$ws.Range("J26").Value = 12000
$ws.Range("L26").Value = 12000
$ws.Range("N26").Value = -12688
$ws.Range("H35").Value = 1925.8572
$ws.Range("H69").Value = 5166.6665
$ws.Range("I69").Value = 4750
$ws.Range("K69").Value = 14250
$ws.Range("M69").Value = -13376
$ws.Range("H72").Value = 5166.6665
$ws.Range("I72").Value = 4750
$ws.Range("K72").Value = 42750
$ws.Range("M72").Value = -38382
$ws.Range("H86").Value = 9000
$ws.Range("I86").Value = 9500
$ws.Range("K86").Value = 9500
$ws.Range("M86").Value = -8377
$ws.Range("H89").Value = 9000
$ws.Range("I89").Value = 9500
$ws.Range("K89").Value = 47500
$ws.Range("M89").Value = -41884
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H137").Value = 2778.3125
$ws.Range("J137").Value = 4459.2856
$ws.Range("L137").Value = 13377.8568
$ws.Range("N137").Value = -18477.8568
$ws.Range("H138").Value = 10030.7705
$ws.Range("I138").Value = 7986.778
$ws.Range("J138").Value = 10502.462
$ws.Range("K138").Value = 23960.334
$ws.Range("L138").Value = 31507.386
$ws.Range("M138").Value = -18820.334
$ws.Range("N138").Value = -41787.386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2291.3333
$ws.Range("I45").Value = 1655.4286
$ws.Range("K45").Value = 1655.4286
$ws.Range("M45").Value = -1278.4286
$ws.Range("H105").Value = 55000
$ws.Range("J105").Value = 55000
$ws.Range("L105").Value = 55000
$ws.Range("N105").Value = -61988
$ws.Range("H110").Value = 3186.111
$ws.Range("I110").Value = 1782.6666
$ws.Range("J110").Value = 5993
$ws.Range("K110").Value = 1782.6666
$ws.Range("L110").Value = 5993
$ws.Range("M110").Value = 262.3334
$ws.Range("N110").Value = -10083
$ws.Range("H132").Value = 1983.4546
$ws.Range("I132").Value = 1581.8
$ws.Range("K132").Value = 4745.4
$ws.Range("M132").Value = -2215.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 750.1667
$ws.Range("I22").Value = 800.2
$ws.Range("K22").Value = 800.2
$ws.Range("M22").Value = -627.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 8750
$ws.Range("J26").Value = 8750
$ws.Range("L26").Value = 8750
$ws.Range("N26").Value = -9324
$ws.Range("H31").Value = 3095.4614
$ws.Range("I31").Value = 2978.5833
$ws.Range("K31").Value = 2978.5833
$ws.Range("M31").Value = -2683.5833
$ws.Range("H34").Value = 3095.4614
$ws.Range("I34").Value = 2978.5833
$ws.Range("K34").Value = 2978.5833
$ws.Range("M34").Value = -2776.5833
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H62").Value = 3981.6667
$ws.Range("J62").Value = 3497.5
$ws.Range("L62").Value = 3497.5
$ws.Range("N62").Value = -4745.5
$ws.Range("H65").Value = 3981.6667
$ws.Range("J65").Value = 3497.5
$ws.Range("L65").Value = 17487.5
$ws.Range("N65").Value = -23727.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7250
$ws.Range("I3").Value = 7250
$ws.Range("K3").Value = 21750
$ws.Range("M3").Value = -21638
$ws.Range("H4").Value = 238
$ws.Range("I4").Value = 142.75
$ws.Range("K4").Value = 428.25
$ws.Range("M4").Value = -316.25
$ws.Range("H21").Value = 1002
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 290
$ws.Range("J23").Value = 350
$ws.Range("L23").Value = 1050
$ws.Range("N23").Value = -1520
$ws.Range("H34").Value = 2846.6667
$ws.Range("J34").Value = 3400
$ws.Range("L34").Value = 10200
$ws.Range("N34").Value = -10368
$ws.Range("H39").Value = 4001.5
$ws.Range("J39").Value = 4001.5
$ws.Range("L39").Value = 12004.5
$ws.Range("N39").Value = -12592.5
$ws.Range("H55").Value = 3125
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 4000
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 12000
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = -12354
$ws.Range("H108").Value = 799.5
$ws.Range("I108").Value = 799.5
$ws.Range("K108").Value = 2398.5
$ws.Range("M108").Value = 481.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 18999.666
$ws.Range("J33").Value = 18999.666
$ws.Range("L33").Value = 18999.666
$ws.Range("N33").Value = -19503.666
$ws.Range("H53").Value = 39
$ws.Range("I53").Value = 39
$ws.Range("K53").Value = 39
$ws.Range("M53").Value = 592
$ws.Range("H80").Value = 3421
$ws.Range("J80").Value = 2949
$ws.Range("L80").Value = 2949
$ws.Range("N80").Value = -4945
$ws.Range("H83").Value = 3421
$ws.Range("J83").Value = 2949
$ws.Range("L83").Value = 14745
$ws.Range("N83").Value = -24729
$ws.Range("H107").Value = 766.6667
$ws.Range("I107").Value = 766.6667
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 766.6667
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1153.3333
$ws.Range("N107").ClearContents()
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 7000
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 1500
$ws.Range("J4").Value = 1500
$ws.Range("L4").Value = 1500
$ws.Range("N4").Value = -1726
$ws.Range("H28").Value = 1500
$ws.Range("J28").Value = 1500
$ws.Range("L28").Value = 1500
$ws.Range("N28").Value = -1964
$ws.Range("H32").Value = 1950
$ws.Range("I32").Value = 1950
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1950
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1633
$ws.Range("N32").ClearContents()
$ws.Range("H37").Value = 1500
$ws.Range("J37").Value = 1500
$ws.Range("L37").Value = 1500
$ws.Range("N37").Value = -1714
$ws.Range("H93").Value = 847
$ws.Range("I93").Value = 847
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 847
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 401
$ws.Range("N93").ClearContents()
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H122").Value = 8996
$ws.Range("I122").Value = 8996
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 26988
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -24538
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 7199.727
$ws.Range("I132").Value = 6024.75
$ws.Range("K132").Value = 18074.25
$ws.Range("M132").Value = -15544.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120
